$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.645.31"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "2.112.45"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.18"
$ws.Range("E5").Value = "  +1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.666"
$ws.Range("E6").Value = "  -5.44%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.27"
$ws.Range("E8").Value = "  +6.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "60.53"
$ws.Range("E9").Value = "  +3.22%  "
$ws.Range("E10").Value = "  +1.13%  "
$ws.Range("E11").Value = "  -2.56%  "
$ws.Range("E12").Value = "  +0.09%  "
$ws.Range("D13").Value = "2.420.17"
$ws.Range("E13").Value = "  +9.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.26"
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.833"
$ws.Range("E15").Value = "  +3.00%  "
$ws.Range("D16").Value = "2.114.74"
$ws.Range("E16").Value = "  +9.86%  "
$ws.Range("E17").Value = "  -0.33%  "
$ws.Range("D18").Value = "36.589.19"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.55"
$ws.Range("E19").Value = "  -1.20%  "
$ws.Range("E20").Value = "  -2.94%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "240.62"
$ws.Range("E22").Value = "  -4.39%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  -6.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.03"
$ws.Range("E26").Value = "  +2.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.59"
$ws.Range("E27").Value = "  +14.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.19"
$ws.Range("E28").Value = "  +4.36%  "
$ws.Range("E29").Value = "  -8.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.21"
$ws.Range("E30").Value = "  +64.55%  "
$ws.Range("E31").Value = "  -4.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.49"
$ws.Range("E32").Value = "  -1.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0942"
$ws.Range("E33").Value = "  +7.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0599"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.956"
$ws.Range("E35").Value = "  +6.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.90"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  +15.01%  "
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.12"
$ws.Range("E39").Value = "  -5.83%  "
$ws.Range("E40").Value = "  -12.15%  "
$ws.Range("E41").Value = "  +6.44%  "
$ws.Range("E42").Value = "  -1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "98.72"
$ws.Range("E43").Value = "  -7.44%  "
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "15.95"
$ws.Range("E45").Value = "  -8.53%  "
$ws.Range("D46").Value = "1.349.63"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.19"
$ws.Range("E47").Value = "  +11.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0842"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").Value = "2.300.01"
$ws.Range("E49").Value = "  +9.41%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("E51").Value = "  -4.61%  "
